# Auto-generated edit script for cryptos.xlsx crypto price/volume update
# Commit: Updated cryptos list on Fri Sep 22 19:50:42 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.630.97"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "1.597.25"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'211.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("E6").Value = "  +0.94%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.244"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.36%  "
$ws.Range("D10").Value = "'19.44"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.18%  "
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").Value = "1.822.69"
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("D13").Value = "1.569.67"
$ws.Range("E13").Value = "  -1.51%  "
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("E15").Value = "  -1.04%  "
$ws.Range("D16").Value = "'64.70"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "26.636.83"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").Value = "'207.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.44%  "
$ws.Range("D21").Value = "'6.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.61%  "
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "'2.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.68%  "
$ws.Range("D24").Value = "'8.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.58%  "
$ws.Range("D25").Value = "'145.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.79%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "'7.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.71%  "
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("D29").Value = "'15.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("E30").Value = "  -0.75%  "
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("E32").Value = "  -0.33%  "
$ws.Range("E33").Value = "  -1.67%  "
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("D35").Value = "1.280.89"
$ws.Range("E35").Value = "  -3.02%  "
$ws.Range("D36").Value = "'2.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.70%  "
$ws.Range("E37").Value = "  -0.43%  "
$ws.Range("E38").Value = "  -0.61%  "
$ws.Range("D39").Value = "'0.842"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.26%  "
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("E42").Value = "  +1.53%  "
$ws.Range("D43").Value = "'0.785"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.46%  "
$ws.Range("D44").Value = "'63.86"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.11%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'0.917"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.36%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.734.18"
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("D47").Value = "'89.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.61%  "
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("E49").Value = "  -1.16%  "
$ws.Range("E50").Value = "  +3.42%  "
$ws.Range("E51").Value = "  -1.43%  "
